# GPLIM-3541: add Material Type as required header for Manifest uploads
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 "Material Type" - copy formatting from the other
# header cell styled the same way (A1, "Specimen_Number").
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "Material Type"

# New data column G2:G24 "DNA:Genomic" - copy formatting from the
# centered data column C (SAMPLE_TYPE values use the same style).
$ws.Range("C2").Copy()
$ws.Range("G2:G24").PasteSpecial(-4122)   # xlPasteFormats
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 7).Value = "DNA:Genomic"
}

$ws.Range("H3").Select()
